$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column headers for columns G (4th sectioning type) and H (5th sectioning type)
$ws.Range("H1").Value = "Frozen sectioning & H&E stain"
$ws.Range("G1").Value = "Frozen sectioning-unstained slide"

# Widen column G to fit the new, longer header text
$ws.Columns.Item(7).ColumnWidth = 35.6

# Update the selected cell in the sheet view
$ws.Range("G9").Select()
